$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")

# Update Date value
$wsMeta.Range("B8").Value = "2025-07-11T12:29:53+00:00"

# Add FR Jurisdiction value
$wsMeta.Range("B11").Value = "FRANCE"
